$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "279.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "6.15%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "1.35%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "4.900"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "4.39%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.06341"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "3.88%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "6.945"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "3.62%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "3.354"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "6.07%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.8820"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "3.82%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.9490"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "4.67%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.1471"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "4.51%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.05110"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "1.69%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.07485"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "5.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.03137"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "0.50%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.09063"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "0.18%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "1.38%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.0006284"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "1.73%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.005754"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "-4.01%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "3.480"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "0.90%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "2.295"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "5.93%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.1309"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "2.29%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "3.861"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "-5.26%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.04320"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "1.77%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.001175"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "-0.43%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.003621"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "-10.76%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "-0.05%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.0001692"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "-12.66%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.04047"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "2.54%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.006618"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "58.15%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "4.75%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.002348"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "11.35%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "7.32%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.00005209"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "2.23%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "0.02%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.377"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "820.89%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "5.94%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "0.02%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "-0.05%"
